# Chale Proposal 2 - update task-tracking worksheet
#
# - Column I ("Ready Now?") for the first three tasks (rows 5-7) was
#   "Python only" and is now simply "Yes" (those strings are no longer
#   referenced anywhere else, so Excel drops them from the shared-strings
#   table on save).
# - E7 ("Status" for the "Report performance" task) used to hold a
#   free-text explanation; it now holds the same kind of numeric/percent
#   value (1 = 100%) used by the sibling rows.
# - The active selection on the sheet moved from F6 to E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "Yes"
$ws.Range("I6").Value = "Yes"
$ws.Range("I7").Value = "Yes"

$ws.Range("E7").Value = 1

$ws.Range("E7").Select()
